# Append the 2025-03-17 price row (row 16) to every Solar_Prices sheet,
# carrying forward each sheet's last known value (row 15, column B).

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "N-Dense",
    "N-Type",
    "N-type Wafer",
    "Cell Topcon 183mm",
    "Module Topcon 183mm",
    "Silver Rear_side",
    "Silver Busbar front-side",
    "Silver finger front-side",
    "USD_CNY"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Last known price (row 15, column B) - carried forward as the new row's value.
    $lastValue = $ws.Range("B15").Value2

    # New date cell (A16). Force text so Excel doesn't auto-convert the
    # "yyyy-mm-dd" string into a date serial number, matching the existing
    # date cells in column A (all stored as literal text).
    $dateCell = $ws.Range("A16")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2025-03-17"
    $dateCell.ClearFormats()

    # New price cell (B16), carried forward from row 15. Force text so a
    # numeric-looking value (e.g. "40") isn't coerced into a real number,
    # matching the existing price cells which are stored as literal text.
    $priceCell = $ws.Range("B16")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $lastValue
    $priceCell.ClearFormats()
}
